$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply yellow-fill styling to the first PID table's D column (rows 13-16) ---
# D13, D14, D15 -> fill only
$ws.Range("D13").Interior.Color = 65535
$ws.Range("D14").Interior.Color = 65535
$ws.Range("D15").Interior.Color = 65535
# D16 -> fill + scientific number format (matches existing "0.00E+00"-style column)
$ws.Range("D16").Interior.Color = 65535
$ws.Range("D16").NumberFormat = "0.00E+00"

# --- Relocate the old bottom summary block (B19:D23) out to S27:U31 ---
$ws.Range("S27").Value = $ws.Range("B19").Value2
$ws.Range("S28").Value = $ws.Range("B20").Value2
$ws.Range("T28").Value = $ws.Range("C20").Value2
$ws.Range("U28").Value = $ws.Range("D20").Value2
$ws.Range("S29").Value = $ws.Range("B21").Value2
$ws.Range("T29").Value = $ws.Range("C21").Value2
$ws.Range("U29").Value = $ws.Range("D21").Value2
$ws.Range("S30").Value = $ws.Range("B22").Value2
$ws.Range("T30").Value = $ws.Range("C22").Value2
$ws.Range("U30").Value = $ws.Range("D22").Value2
$ws.Range("S31").Value = $ws.Range("B23").Value2
$ws.Range("T31").Value = $ws.Range("C23").Value2
$ws.Range("U31").Value = $ws.Range("D23").Value2
$ws.Range("T31").NumberFormat = "0.00E+00"
$ws.Range("U31").NumberFormat = "0.00E+00"

# Wipe the old B19:D23 block completely (values + formats); column E (E23) is left as-is
$ws.Range("B19:D23").Clear()

# --- Insert the new "scanRate / size" dated block at A21:C26 ---
$ws.Range("A21").Value = 20160902

# Shared-string insertion order matters for an exact match: "scanRate" (row23)
# was authored before "size"/"30um,15um"/"15um" (row22), so write row 23 first.
$ws.Range("A23").Value = "scanRate"
$ws.Range("B23").Value = 0.2
$ws.Range("C23").Value = 0.5

$ws.Range("A22").Value = "size"
$ws.Range("B22").Value = "30um,15um"
$ws.Range("C22").Value = "15um"

$ws.Range("A24").Value = "P"
$ws.Range("B24").Value = 0.4
$ws.Range("B24").Interior.Color = 65535

$ws.Range("A25").Value = "I"
$ws.Range("B25").Value = 0.3
$ws.Range("B25").Interior.Color = 65535

$ws.Range("A26").Value = "D"
$ws.Range("B26").Value = 0.00005
$ws.Range("B26").Interior.Color = 65535
$ws.Range("B26").NumberFormat = "0.00E+00"

# --- Selection + print setup to match the final saved state ---
[void]$ws.Range("C24").Select()
$ws.PageSetup.Orientation = 1
